$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 600
$ws.Range("I2").Value = 800
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -687
$ws.Range("N2").Value = -626
$ws.Range("H4").Value = 8358.200000000001
$ws.Range("I4").Value = 7948
$ws.Range("K4").Value = 7948
$ws.Range("M4").Value = -7834
$ws.Range("H9").Value = 763.7273
$ws.Range("I9").Value = 740.1
$ws.Range("K9").Value = 740.1
$ws.Range("M9").Value = -571.1
$ws.Range("H17").Value = 1199
$ws.Range("J17").Value = 1199
$ws.Range("L17").Value = 3597
$ws.Range("N17").Value = -3933
$ws.Range("H33").Value = 1431.1
$ws.Range("I33").Value = 1200.5714
$ws.Range("J33").Value = 1969
$ws.Range("K33").Value = 1200.5714
$ws.Range("L33").Value = 1969
$ws.Range("M33").Value = -971.5714
$ws.Range("N33").Value = -2427
$ws.Range("H53").Value = 450.125
$ws.Range("J53").Value = 493
$ws.Range("L53").Value = 493
$ws.Range("N53").Value = -1767
$ws.Range("H69").Value = 22218.334
$ws.Range("I69").Value = 11996.167
$ws.Range("K69").Value = 35988.501
$ws.Range("M69").Value = -35114.501
$ws.Range("H72").Value = 22218.334
$ws.Range("I72").Value = 11996.167
$ws.Range("K72").Value = 107965.503
$ws.Range("M72").Value = -103597.503
$ws.Range("H87").Value = 126950
$ws.Range("J87").Value = 126950
$ws.Range("L87").Value = 126950
$ws.Range("N87").Value = -129446
$ws.Range("H90").Value = 126950
$ws.Range("J90").Value = 126950
$ws.Range("L90").Value = 380850
$ws.Range("N90").Value = -393330
$ws.Range("H103").Value = 894.5714
$ws.Range("I103").Value = 502.1111
$ws.Range("K103").Value = 1506.3333
$ws.Range("M103").Value = -920.3333
$ws.Range("H127").Value = 1507.6666
$ws.Range("I127").Value = 1224.1428
$ws.Range("K127").Value = 3672.4284
$ws.Range("M127").Value = 1287.5716
$ws.Range("H135").Value = 1038.3103
$ws.Range("I135").Value = 1041.1482
$ws.Range("K135").Value = 9370.3338
$ws.Range("M135").Value = -6835.3338
$ws.Range("H138").Value = 2339.1538
$ws.Range("J138").Value = 2605.3809
$ws.Range("L138").Value = 7816.1427
$ws.Range("N138").Value = -18096.1427

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 80004.5
$ws.Range("J8").Value = 80004.5
$ws.Range("L8").Value = 80004.5
$ws.Range("N8").Value = -80292.5
$ws.Range("H32").Value = 10640532
$ws.Range("I32").Value = 10640532
$ws.Range("K32").Value = 10640532
$ws.Range("M32").Value = -10640245
$ws.Range("H63").Value = 5214.6665
$ws.Range("I63").Value = 4972
$ws.Range("K63").Value = 4972
$ws.Range("M63").Value = -4286
$ws.Range("H66").Value = 5214.6665
$ws.Range("I66").Value = 4972
$ws.Range("K66").Value = 24860
$ws.Range("M66").Value = -21428
$ws.Range("H74").Value = 9623480
$ws.Range("I74").Value = 14707352
$ws.Range("J74").Value = 20610.555
$ws.Range("K74").Value = 14707352
$ws.Range("L74").Value = 20610.555
$ws.Range("M74").Value = -14706478
$ws.Range("N74").Value = -22358.555
$ws.Range("H77").Value = 9623480
$ws.Range("I77").Value = 14707352
$ws.Range("J77").Value = 20610.555
$ws.Range("K77").Value = 73536760
$ws.Range("L77").Value = 103052.775
$ws.Range("M77").Value = -73532392
$ws.Range("N77").Value = -111788.775
$ws.Range("H97").Value = 1095.5769
$ws.Range("I97").Value = 1095.5769
$ws.Range("K97").Value = 1095.5769
$ws.Range("M97").Value = -599.5769
$ws.Range("H110").Value = 1320.8334
$ws.Range("I110").Value = 1425
$ws.Range("K110").Value = 1425
$ws.Range("M110").Value = 620
$ws.Range("H132").Value = 8333
$ws.Range("I132").Value = 3511.3076
$ws.Range("J132").Value = 18780
$ws.Range("K132").Value = 10533.9228
$ws.Range("L132").Value = 56340
$ws.Range("M132").Value = -8003.9228
$ws.Range("N132").Value = -61400

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5750
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -1387
$ws.Range("N7").Value = -10226
$ws.Range("H27").Value = 64994.5
$ws.Range("J27").Value = 64994.5
$ws.Range("L27").Value = 64994.5
$ws.Range("N27").Value = -65378.5
$ws.Range("H80").Value = 1312.8
$ws.Range("J80").Value = 1362.2727
$ws.Range("L80").Value = 1362.2727
$ws.Range("N80").Value = -3358.2727
$ws.Range("H83").Value = 1312.8
$ws.Range("J83").Value = 1362.2727
$ws.Range("L83").Value = 6811.363499999999
$ws.Range("N83").Value = -16795.3635
$ws.Range("H105").Value = 2658.7
$ws.Range("I105").Value = 2573.375
$ws.Range("K105").Value = 2573.375
$ws.Range("M105").Value = -826.375
$ws.Range("H107").Value = 2574.6667
$ws.Range("I107").Value = 1659.1333
$ws.Range("J107").Value = 7152.3335
$ws.Range("K107").Value = 1659.1333
$ws.Range("L107").Value = 7152.3335
$ws.Range("M107").Value = 260.8667
$ws.Range("N107").Value = -10992.3335

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 64103
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 64103
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 64103
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -64329
$ws.Range("H22").Value = 774.5
$ws.Range("I22").Value = 774.5
$ws.Range("K22").Value = 774.5
$ws.Range("M22").Value = -424.5
$ws.Range("H31").Value = 558655.3
$ws.Range("I31").Value = 3701.8333
$ws.Range("K31").Value = 3701.8333
$ws.Range("M31").Value = -3406.8333
$ws.Range("H34").Value = 558655.3
$ws.Range("I34").Value = 3701.8333
$ws.Range("K34").Value = 3701.8333
$ws.Range("M34").Value = -3499.8333
$ws.Range("H134").Value = 337013.16
$ws.Range("I134").Value = 477079.47
$ws.Range("J134").Value = 10191.777
$ws.Range("K134").Value = 1431238.41
$ws.Range("L134").Value = 30575.331
$ws.Range("M134").Value = -1428703.41
$ws.Range("N134").Value = -35645.331

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4883113.5
$ws.Range("J4").Value = 91690.82000000001
$ws.Range("L4").Value = 275072.46
$ws.Range("N4").Value = -275296.46
$ws.Range("H132").Value = 1275.6364
$ws.Range("J132").Value = 1335.2727
$ws.Range("L132").Value = 12017.4543
$ws.Range("N132").Value = -17077.4543
$ws.Range("H133").Value = 6061.75
$ws.Range("J133").Value = 6923.5
$ws.Range("L133").Value = 20770.5
$ws.Range("N133").Value = -30890.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 35004.5
$ws.Range("I5").Value = 30004
$ws.Range("J5").Value = 40005
$ws.Range("K5").Value = 30004
$ws.Range("L5").Value = 40005
$ws.Range("M5").Value = -29892
$ws.Range("N5").Value = -40229
$ws.Range("H82").Value = 113999.5
$ws.Range("J82").Value = 128000
$ws.Range("L82").Value = 128000
$ws.Range("N82").Value = -128766
$ws.Range("H85").Value = 113999.5
$ws.Range("J85").Value = 128000
$ws.Range("L85").Value = 128000
$ws.Range("N85").Value = -130652
$ws.Range("H97").Value = 1131.5714
$ws.Range("I97").Value = 1277.5333
$ws.Range("K97").Value = 1277.5333
$ws.Range("M97").Value = -781.5333000000001
$ws.Range("H102").Value = 3341.923
$ws.Range("I102").Value = 2548.9092
$ws.Range("J102").Value = 7703.5
$ws.Range("K102").Value = 2548.9092
$ws.Range("L102").Value = 7703.5
$ws.Range("M102").Value = -926.9092000000001
$ws.Range("N102").Value = -10947.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1946
$ws.Range("J22").Value = 299.5
$ws.Range("L22").Value = 299.5
$ws.Range("N22").Value = -889.5
$ws.Range("H27").Value = 1946
$ws.Range("J27").Value = 299.5
$ws.Range("L27").Value = 299.5
$ws.Range("N27").Value = -513.5
$ws.Range("H136").Value = 37295.137
$ws.Range("I136").Value = 5812.846
$ws.Range("K136").Value = 17438.538
$ws.Range("M136").Value = -14888.538

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1265.8
$ws.Range("I113").Value = 1244.1428
$ws.Range("K113").Value = 3732.4284
$ws.Range("M113").Value = -1562.4284
